# Move the last data row (row 41, BNr "Z03_B04") up to row 8,
# shifting the rows currently at 8..40 down by one row (to 9..41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values of row 41 (the row that needs to move up to row 8)
$lastA = $ws.Range("A41").Value2
$lastB = $ws.Range("B41").Value2
$lastC = $ws.Range("C41").Value2
$lastD = $ws.Range("D41").Value2

# Shift rows 8..40 down to 9..41 (work from the bottom up to avoid overwrites)
for ($r = 40; $r -ge 8; $r--) {
    $a = $ws.Range("A" + $r).Value2
    $b = $ws.Range("B" + $r).Value2
    $c = $ws.Range("C" + $r).Value2
    $d = $ws.Range("D" + $r).Value2
    $ws.Range("A" + ($r + 1)).Value = $a
    $ws.Range("B" + ($r + 1)).Value = $b
    $ws.Range("C" + ($r + 1)).Value = $c
    $ws.Range("D" + ($r + 1)).Value = $d
}

# Place the captured last row values into row 8
$ws.Range("A8").Value = $lastA
$ws.Range("B8").Value = $lastB
$ws.Range("C8").Value = $lastC
$ws.Range("D8").Value = $lastD
